$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BS is column 71 (A=1 ... Z=26, AA=27 ... BR=70, BS=71, BT=72, BU=73).
# Insert a brand-new column there; this shifts the existing BS ("nom") to BT
# and the existing BT ("url_produit") to BU, and widens the used range to
# BU206 automatically. It also copies the BS column's old style (s="1" header
# style, etc.) is not something Insert does by itself for data cells, but for
# row 1 it correctly carries the header cell style rightward along with the
# values, and leaves the freshly inserted column empty/blank.
$ws.Columns.Item(71).Insert()

# New header cell for the inserted timestamp column.
$ws.Range("BS1").Value = "2026-01-30 21:18:44"

# For the product rows that still have an active price in BR (rows 2-80),
# carry that last known price forward into the newly inserted BS column so
# the new timestamp snapshot repeats the most recent price, matching the
# other timestamp columns' "carry forward" behaviour.
for ($row = 2; $row -le 80; $row++) {
    $lastPrice = $ws.Cells.Item($row, 70).Value2
    $ws.Cells.Item($row, 71).Value = $lastPrice
}
